# ClientPortal_Data.xlsx maintenance update
# - Customer (sheet3) & Test (sheet4): insert a new "Status" column ahead of
#   "Type", normalise blank "Type" values to "Normal", and repurpose the old
#   "Store" column into a "Class" column (constant "CCB").
# - Test sheet gets refreshed sample rows (Maria Lopez / Emily Johnson)
#   replacing the old single John Doe row, and becomes the active tab.
# - ValueList selection reset to A3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: insert the Status/Type/Class column trio on a Flag/Id/.../Phone
# sheet (Customer & Test share the exact same 19-column layout), then return
# the worksheet so callers can keep chaining.
# ---------------------------------------------------------------------------
function Expand-CustomerColumns($ws, $lastDataRow) {
    # Insert a new column in front of the old "Type" column (column J).
    $ws.Columns.Item(10).Insert()

    # New column J = Status, constant "Active" for every data row.
    $ws.Range("J1").Value = "Status"
    $ws.Range("J2:J$lastDataRow").Value = "Active"

    # Old "Type" data shifted into K; fill the previously blank cells with
    # the new default "Normal".
    $ws.Range("K1").Value = "Type"
    for ($r = 2; $r -le $lastDataRow; $r++) {
        $cell = $ws.Range("K$r")
        if ([string]::IsNullOrEmpty($cell.Value2)) {
            $cell.Value = "Normal"
        }
    }

    # Old "Store" column (now L) becomes "Class" with a constant "CCB".
    $ws.Range("L1").Value = "Class"
    $ws.Range("L2:L$lastDataRow").Value = "CCB"

    # Give the Type column its own left/top/bottom thin border, matching the
    # new Status/Type/Class grouping.
    $typeRange = $ws.Range("K2:K$lastDataRow")
    $typeRange.Borders.Item(7).LineStyle = 1
    $typeRange.Borders.Item(7).Weight = 2
    $typeRange.Borders.Item(8).LineStyle = 1
    $typeRange.Borders.Item(8).Weight = 2
    $typeRange.Borders.Item(9).LineStyle = 1
    $typeRange.Borders.Item(9).Weight = 2

    return $ws
}

# Helper: (re)point every mailto hyperlink on a sheet at the Email column,
# one row per supplied address, in row order starting at row 2.
function Reset-EmailHyperlinks($ws, [string[]]$addresses) {
    $ws.Hyperlinks.Delete()
    $r = 2
    foreach ($addr in $addresses) {
        $ws.Hyperlinks.Add($ws.Range("S$r"), "mailto:$addr")
        $r = $r + 1
    }
}

# ---------------------------------------------------------------------------
# Customer sheet
# ---------------------------------------------------------------------------
$wsCustomer = $wb.Worksheets.Item("Customer")
Expand-CustomerColumns $wsCustomer 7 | Out-Null

# William's email loses its stray trailing tab character.
$wsCustomer.Range("S7").Value = "william@gmail.com"

Reset-EmailHyperlinks $wsCustomer @(
    "john.doe@example.com",
    "maria.lopez@example.es",
    "emily@gmail.com",
    "linda.brown@example.com",
    "mike@hotmail.com",
    "william@gmail.com"
)

$wsCustomer.Range("C12").Select() | Out-Null

# ---------------------------------------------------------------------------
# Test sheet — replace the single John Doe sample row with the Maria
# Lopez / Emily Johnson rows, after the same column expansion.
# ---------------------------------------------------------------------------
$wsTest = $wb.Worksheets.Item("Test")
Expand-CustomerColumns $wsTest 2 | Out-Null

$wsTest.Range("A2").Value = "Yes"
$wsTest.Range("B2").Value = 2
$wsTest.Range("C2").Value = "Maria Lopez"
$wsTest.Range("D2").Value = "Calle Mayor, 15"
$wsTest.Range("E2").Value = "Apt 3A"
$wsTest.Range("F2").Value = "Spain"
$wsTest.Range("G2").Value = "Madrid"
$wsTest.Range("H2").Value = "Madrid"
$wsTest.Range("I2").Value = 28013
$wsTest.Range("J2").Value = "Active"
$wsTest.Range("K2").Value = "Normal"
$wsTest.Range("L2").Value = "CCB"
$wsTest.Range("M2").Value = "No"
$wsTest.Range("N2").Value = "Yes"
$wsTest.Range("O2").Value = 25
$wsTest.Range("P2").Value = 6
$wsTest.Range("Q2").Value = "N/A"
$wsTest.Range("R2").Value = "N/A"
$wsTest.Range("S2").Value = "maria.lopez@example.es"
$wsTest.Range("T2").Value = "+34915556789"

$wsTest.Range("A3").Value = "Yes"
$wsTest.Range("B3").Value = ""
$wsTest.Range("C3").Value = "Emily Johnson"
$wsTest.Range("D3").Value = "600 Congress Avenue"
$wsTest.Range("E3").Value = "Suite 14"
$wsTest.Range("F3").Value = "United States of America"
$wsTest.Range("G3").Value = "Florida"
$wsTest.Range("H3").Value = "Hawthorne"
$wsTest.Range("I3").Value = 78701
$wsTest.Range("J3").Value = "Active"
$wsTest.Range("K3").Value = "Normal"
$wsTest.Range("L3").Value = "CCB"
$wsTest.Range("M3").Value = "Yes"
$wsTest.Range("N3").Value = "Yes"
$wsTest.Range("O3").Value = 25
$wsTest.Range("P3").Value = 0
$wsTest.Range("Q3").Value = "N/A"
$wsTest.Range("R3").Value = "N/A"
$wsTest.Range("S3").Value = "emily@gmail.com"
$wsTest.Range("T3").Value = "+34955559145"

Reset-EmailHyperlinks $wsTest @(
    "maria.lopez@example.es",
    "emily@gmail.com"
)

# ---------------------------------------------------------------------------
# Window / selection bookkeeping
# ---------------------------------------------------------------------------
$wsValueList = $wb.Worksheets.Item("ValueList")
$wsValueList.Range("A3").Select() | Out-Null

# Test becomes the front-most (active) tab; its own selection moves to N3.
$wsTest.Activate() | Out-Null
$wsTest.Range("N3").Select() | Out-Null
